$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column N (14th column); shifts old N -> O.
$ws.Columns("N:N").Insert()

# 2. New column N width + style (matches existing data columns: width 12.140625, style index 2)
$ws.Columns("N:N").ColumnWidth = 12.140625

# 3. New header text in N1 (top header row) - new issue note
$ws.Range("N1").Value = "in main fun, check that intern fun have the good args"

# 4. Append regex note to K1 (existing header cell text)
$ws.Range("K1").Value = "80`n. preceeded by :::`n\b(?!base\b)\w+(?=::)"

# 5. New content cell D10 = "issue 82"
$ws.Range("D10").Value = "issue 82"

# 6. Fill in the rest of row 10 (C10,E10,G10:M10) with "x" to match the extended row
$ws.Range("C10").Value = "x"
$ws.Range("E10").Value = "x"
$ws.Range("G10").Value = "x"
$ws.Range("H10").Value = "x"
$ws.Range("I10").Value = "x"
$ws.Range("J10").Value = "x"
$ws.Range("K10").Value = "x"
$ws.Range("L10").Value = "x"
$ws.Range("M10").Value = "x"

# 7. Fill in the rest of row 17 (C17,E17,G17:M17) with "x"
$ws.Range("C17").Value = "x"
$ws.Range("E17").Value = "x"
$ws.Range("G17").Value = "x"
$ws.Range("H17").Value = "x"
$ws.Range("I17").Value = "x"
$ws.Range("J17").Value = "x"
$ws.Range("K17").Value = "x"
$ws.Range("L17").Value = "x"
$ws.Range("M17").Value = "x"

# 8. Highlight (yellow fill) the function-name column (A) for the "not yet done" rows
$yellow = 65535
$highlightRows = @(2,3,4,5,6,7,19,20,21,22)
foreach ($r in $highlightRows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Interior.Color = $yellow
}

# 9. Highlight C17 with yellow fill + centered alignment (style 4)
$c17 = $ws.Range("C17")
$c17.Interior.Color = $yellow
$c17.HorizontalAlignment = -4108
$c17.VerticalAlignment = -4108
$c17.WrapText = $true

# 10. Update the selection to N18 (matches the saved cursor position in the diff)
$ws.Range("N18").Select()
